# Add the I0 / IF columns to the save-data sheet.
#
# H1 ("IP") is the last existing header cell, styled with the bold/bordered
# header style (style index 1). The new I1/J1 header cells need that same
# style, and I2/J2 are plain numeric data cells (no style), matching H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H1's formatting (only) onto the new header cells so they pick up the
# existing header style instead of minting a brand new one.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values (row 2), unstyled like the other numeric cells in row 2.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
